$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column (C) from row 2 through row 292
# The value stored is an Excel serial date number; it is being bumped by 1 day.
$ws.Range("C2:C292").Value = 45172
